$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (A1) onto the three new header cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record values for every data row (2-38)
for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 30).Value = 89
    $ws.Cells.Item($row, 31).Value = 73
    $ws.Cells.Item($row, 32).Value = 1
}
